$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header-like numeric labels for combination columns
$ws.Range("H1").Value = 1
$ws.Range("I1").Value = 2
$ws.Range("J1").Value = 3
$ws.Range("K1").Value = 12
$ws.Range("L1").Value = 13
$ws.Range("M1").Value = 23
$ws.Range("N1").Value = 123

# Row 2: formulas (not part of the shared-formula group)
$ws.Range("H2").Formula = "=A2"
$ws.Range("I2").Formula = "=C2"
$ws.Range("J2").Formula = "=E2"
$ws.Range("K2").Formula = "=A2*C2"
$ws.Range("L2").Formula = "=A2*E2"
$ws.Range("M2").Formula = "=C2*E2"
$ws.Range("N2").Formula = "=A2*C2*E2"

# Rows 3-12: formulas (shared-formula group)
for ($r = 3; $r -le 12; $r++) {
    $ws.Range("H$r").Formula = "=A$r"
    $ws.Range("I$r").Formula = "=C$r"
    $ws.Range("J$r").Formula = "=E$r"
    $ws.Range("K$r").Formula = "=A$r*C$r"
    $ws.Range("L$r").Formula = "=A$r*E$r"
    $ws.Range("M$r").Formula = "=C$r*E$r"
    $ws.Range("N$r").Formula = "=A$r*C$r*E$r"
}

[void]$ws.Range("L13").Select()
